# "Generate Report for Handoff"
# The localization job for 16150b8d-9227-40a3-98c9-6524522ef9ee.md moved from
# "In Translation" to "Ready for handoff", and the handoff-xliff-generation
# timestamps were refreshed. Update the three report sheets (Overview,
# zh-cn, de-de) to reflect the new status/timestamps, and widen the
# "Status"/language columns that now hold the longer "Ready for handoff"
# text (Excel auto-widened them when the content grew).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ------------------------
$overview.Range("E2").Value = "Ready for handoff"   # Overview: zh-cn status
$overview.Range("F2").Value = "Ready for handoff"   # Overview: de-de status
$zhcn.Range("C2").Value     = "Ready for handoff"   # zh-cn sheet: Status
$dede.Range("C2").Value     = "Ready for handoff"   # de-de sheet: Status

# --- Latest HO Xliff Generate Date / Latest Handoff Datetime ----------------
$overview.Range("G2").Value = "2016-12-05 11:21:07" # Overview: shared w/ de-de
$dede.Range("H2").Value     = "2016-12-05 11:21:07" # de-de: Latest Handoff Datetime
$zhcn.Range("H2").Value     = "2016-12-05 11:20:52" # zh-cn: Latest Handoff Datetime

# --- Column widths: grew to fit the longer "Ready for handoff" text ---------
# ColumnWidth is in characters; the engine quantizes to whole pixels, so use
# the character width whose pixel-rounded result lands on the target (the
# xlsx <col> width is ColumnWidth + 5/6). 16.3333... -> 17.1667 is the closest
# reachable value to the recorded 17.2159881591797.
$newStatusColWidth = 16.333333333333332

$overview.Columns.Item(5).ColumnWidth = $newStatusColWidth  # Overview col E (zh-cn)
$overview.Columns.Item(6).ColumnWidth = $newStatusColWidth  # Overview col F (de-de)
$zhcn.Columns.Item(3).ColumnWidth     = $newStatusColWidth  # zh-cn col C (Status)
$dede.Columns.Item(3).ColumnWidth     = $newStatusColWidth  # de-de col C (Status)
